$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = 'Bonez'
$ws.Range('F2').Value = 'Padge'
$ws.Range('G2').Value = 'Chicken'
$ws.Range('H2').Value = 'Nom Nom'
$ws.Range('I2').Value = 'Nodder'
$ws.Range('E3').Value = 'Padge'
$ws.Range('F3').Value = 'Nodder'
$ws.Range('G3').Value = 'Nono'
$ws.Range('H3').Value = 'Bonez'
$ws.Range('I3').Value = 'Padge'
$ws.Range('E4').Value = 'Hoops'
$ws.Range('I4').Value = 'Pizza'
$ws.Range('E5').Value = 'GoGo'
$ws.Range('I5').Value = 'Zombie'
$ws.Range('E6').Value = 'Gaggles'
$ws.Range('E7').Value = 'Opps'
$ws.Range('F7').Value = 'Smiles'
$ws.Range('G7').Value = 'T-Whisk'
$ws.Range('H7').Value = 'Smiles'
$ws.Range('I9').Value = 'Laddy'
$ws.Range('I10').Value = 'Blister'
$ws.Range('I11').Value = 'Bonez'
$ws.Range('I12').Value = 'Sea'
$ws.Range('F13').Value = 'Tross'
$ws.Range('F14').Value = 'Bonez'
$ws.Range('E15').Value = 'Surf Co'
$ws.Range('F15').Value = 'Sea'
$ws.Range('E16').Value = 'Sea'
$ws.Range('F16').Value = 'Indi'
$ws.Range('F17').Value = 'Smiles'
$ws.Range('H17').Value = 'Laddy'
$ws.Range('I17').Value = 'Tross'
$ws.Range('E18').Value = 'Jaws'
$ws.Range('F18').Value = 'Bow Wow'
$ws.Range('H18').Value = 'Burning Bush'
$ws.Range('I18').Value = 'Indi'
$ws.Range('F20').Value = 'Surf Co'
$ws.Range('H22').Value = 'Surf Co'
$ws.Range('H23').Value = 'Bonez'
$ws.Range('F25').Value = 'Sea'
$ws.Range('G25').Value = 'Tross'
$ws.Range('H25').Value = 'Indi'
$ws.Range('D26').Value = 'Goobie'
$ws.Range('E26').Value = 'Hoops'
$ws.Range('F26').Value = 'Unicorn'
$ws.Range('G26').Value = 'Sea'
$ws.Range('D27').Value = 'Bonez'
$ws.Range('E27').Value = 'Ups'
$ws.Range('F27').Value = 'Goobie'
$ws.Range('G27').Value = 'Goobie'
$ws.Range('I27').Value = 'Sea'
$ws.Range('D28').Value = 'Sea'
$ws.Range('E28').Value = 'GoGo'
$ws.Range('F28').Value = 'Zombie'
$ws.Range('G28').Value = 'Zombie'
$ws.Range('I28').Value = 'Zombie'
$ws.Range('F29').Value = 'Tross'
$ws.Range('G29').Value = 'Hoops'
$ws.Range('H29').Value = 'Surf Co'
$ws.Range('I29').Value = 'Hoops'
$ws.Range('E30').Value = 'Nom Nom'
$ws.Range('F30').Value = 'Indi'
$ws.Range('G30').Value = 'Ups'
$ws.Range('H30').Value = 'Bonez'
$ws.Range('I30').Value = 'Ups'
$ws.Range('F31').Value = 'G-Poppy'
$ws.Range('H31').Value = 'Bow Wow'
$ws.Range('F32').Value = 'T-Whisk'
$ws.Range('H32').Value = 'Unicorn'
$ws.Range('E33').Value = 'Smiles'
$ws.Range('F33').Value = 'Bow Wow'
$ws.Range('I33').Value = 'Bow Wow'
$ws.Range('E34').Value = 'Chicken'
$ws.Range('F34').Value = 'Chicken'
$ws.Range('H34').Value = 'T-Whisk'
$ws.Range('E35').Value = 'NONE FOUND'
$ws.Range('F35').Value = 'NONE FOUND'
$ws.Range('H35').Value = 'Chicken'
$ws.Range('I35').Value = 'Smiles'
$ws.Range('E36').Value = 'Pizza'
$ws.Range('G36').Value = 'Pizza'
$ws.Range('D37').Value = 'Blister'
$ws.Range('E37').Value = 'Captain'
$ws.Range('F37').Value = 'Padge'
$ws.Range('G37').Value = 'G-Poppy'
$ws.Range('H37').Value = 'Padge'
$ws.Range('I37').Value = 'Captain'
$ws.Range('E38').Value = 'Nodder'
$ws.Range('I38').Value = 'Hawma'
$ws.Range('E39').Value = 'Stastro'
$ws.Range('I39').Value = 'Pizza'
$ws.Range('E40').Value = 'Tross'
$ws.Range('G40').Value = 'Opps'
$ws.Range('I40').Value = 'Padge'
$ws.Range('E41').Value = 'Hawma'
$ws.Range('F41').Value = 'Burning Bush'
$ws.Range('H41').Value = 'Blister'
$ws.Range('I41').Value = 'Stastro'
$ws.Range('H42').Value = 'Nodder'
$ws.Range('I42').Value = 'Blister'
$ws.Range('E43').Value = 'Laddy'
$ws.Range('F43').Value = 'Blister'
$ws.Range('G43').Value = 'NONE FOUND'
$ws.Range('H43').Value = 'Stastro'
$ws.Range('I43').Value = 'Nono'
$ws.Range('E44').Value = 'Burning Bush'
$ws.Range('F44').Value = 'NONE FOUND'
$ws.Range('G44').Value = 'NONE FOUND'
$ws.Range('H44').Value = 'Nono'
$ws.Range('I44').Value = 'G-Poppy'
$ws.Range('E45').Value = 'NONE FOUND'
$ws.Range('H45').Value = 'Pizza'
$ws.Range('I45').Value = 'Opps'
$ws.Range('E46').Value = 'NONE FOUND'
$ws.Range('I46').Value = 'Laddy'
$ws.Range('G47').Value = 'Unicorn'
$ws.Range('I47').Value = 'GoGo'
$ws.Range('F48').Value = 'Jaws'
$ws.Range('H48').Value = 'Laddy'
$ws.Range('H49').Value = 'Ups'
$ws.Range('I49').Value = 'Nom Nom'
$ws.Range('F50').Value = 'Captain'
$ws.Range('H50').Value = 'Tross'
$ws.Range('F51').Value = 'NONE FOUND'
$ws.Range('H51').Value = 'Jaws'
$ws.Range('E52').Value = 'NONE FOUND'
$ws.Range('G52').Value = 'NONE FOUND'
$ws.Range('F54').Value = 'Surf Co'
$ws.Range('F55').Value = 'Goobie'
$ws.Range('D56').Value = 'Unicorn'
$ws.Range('E56').Value = 'Surf Co'
$ws.Range('D57').Value = 'Surf Co'
$ws.Range('E57').Value = 'Goobie'
$ws.Range('D58').Value = 'Bonez'
$ws.Range('E58').Value = 'Zombie'
$ws.Range('E59').Value = 'Nodder'
$ws.Range('F59').Value = 'Unicorn'
$ws.Range('H59').Value = 'Bow Wow'
$ws.Range('D60').Value = 'Smiles'
$ws.Range('E60').Value = 'Smiles'
$ws.Range('F60').Value = 'Ups'
$ws.Range('G60').Value = 'Nom Nom'
$ws.Range('H60').Value = 'Smiles'
$ws.Range('I60').Value = 'Surf Co'
$ws.Range('D62').Value = 'Nono'
$ws.Range('E62').Value = 'Captain'
$ws.Range('F62').Value = 'Nodder'
$ws.Range('G62').Value = 'Blister'
$ws.Range('H62').Value = 'Jaws'
$ws.Range('I62').Value = 'Padge'
